# Update the recomputed slope-stability slice values (columns AH, AI, AJ,
# BA, BB, BG, BH) for rows 2-22 on the active worksheet, matching a re-run
# of the "spencer" right-facing-slope debugging calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AH2").Value = -1238.89293509685
$ws.Range("AJ2").Value = -12.52138795235413
$ws.Range("BA2").Value = 2.525643448558156
$ws.Range("BB2").Value = -1306.467707567553
$ws.Range("BH2").Value = 81.10397589356462

$ws.Range("AH3").Value = 9178.426058408142
$ws.Range("AI3").Value = 1306.467707567553
$ws.Range("AJ3").Value = -12.52138795235413
$ws.Range("BA3").Value = 1.97470111438668
$ws.Range("BB3").Value = 14619.35178207635
$ws.Range("BG3").Value = 81.10397589356462
$ws.Range("BH3").Value = 68.87335299518907

$ws.Range("AH4").Value = 50902.93712701974
$ws.Range("AI4").Value = -13312.8840745088
$ws.Range("AJ4").Value = -12.52138795235413
$ws.Range("BA4").Value = 1.287292231038373
$ws.Range("BB4").Value = 72130.61487837438
$ws.Range("BG4").Value = 68.87335299518907
$ws.Range("BH4").Value = 51.29604027399181

$ws.Range("AH5").Value = 40737.30859266572
$ws.Range("AI5").Value = -85443.49895288318
$ws.Range("AJ5").Value = -12.52138795235413
$ws.Range("BA5").Value = 1.19620909154821
$ws.Range("BB5").Value = 45999.73653439412
$ws.Range("BG5").Value = 51.29604027399181
$ws.Range("BH5").Value = 43.67869361320182

$ws.Range("AH6").Value = 44147.78855750636
$ws.Range("AI6").Value = -131443.2354872773
$ws.Range("AJ6").Value = -12.52138795235413
$ws.Range("BA6").Value = 1.08969729518768
$ws.Range("BB6").Value = 39649.88084506032
$ws.Range("BG6").Value = 43.67869361320182
$ws.Range("BH6").Value = 37.57150888578319

$ws.Range("AH7").Value = 46724.66041240805
$ws.Range("AI7").Value = -171093.1163323376
$ws.Range("AJ7").Value = -12.52138795235413
$ws.Range("BA7").Value = 1.022798475221807
$ws.Range("BB7").Value = 31435.85858423841
$ws.Range("BG7").Value = 37.57150888578319
$ws.Range("BH7").Value = 32.63926217169369

$ws.Range("AH8").Value = 48888.13030562105
$ws.Range("AI8").Value = -202528.974916576
$ws.Range("AJ8").Value = -12.52138795235413
$ws.Range("BA8").Value = 0.9790507117737064
$ws.Range("BB8").Value = 22422.59851412284
$ws.Range("BG8").Value = 32.63926217169369
$ws.Range("BH8").Value = 28.67038553945456

$ws.Range("AH9").Value = 50700.3740035157
$ws.Range("AI9").Value = -224951.5734306988
$ws.Range("AJ9").Value = -12.52138795235413
$ws.Range("BA9").Value = 0.9507818989724603
$ws.Range("BB9").Value = 13101.04485691148
$ws.Range("BG9").Value = 28.67038553945456
$ws.Range("BH9").Value = 25.5428297092451

$ws.Range("AH10").Value = 52534.46950646966
$ws.Range("AI10").Value = -238052.6182876103
$ws.Range("AJ10").Value = -12.52138795235413
$ws.Range("BA10").Value = 0.9340007043994232
$ws.Range("BB10").Value = 3538.823694418776
$ws.Range("BG10").Value = 25.5428297092451
$ws.Range("BH10").Value = 23.21346531281761

$ws.Range("AH11").Value = 45571.09884572247
$ws.Range("AI11").Value = -241591.4419820291
$ws.Range("AJ11").Value = -12.52138795235413
$ws.Range("BA11").Value = 0.9267320521775275
$ws.Range("BB11").Value = -4431.075802492084
$ws.Range("BG11").Value = 23.21346531281761
$ws.Range("BH11").Value = 21.76039464570018

$ws.Range("AH12").Value = 40617.54697997558
$ws.Range("AI12").Value = -237160.366179537
$ws.Range("AJ12").Value = -12.52138795235413
$ws.Range("BA12").Value = 0.9267441586507567
$ws.Range("BB12").Value = -10079.08169338501
$ws.Range("BG12").Value = 21.76039464570018
$ws.Range("BH12").Value = 20.82149048541144

$ws.Range("AH13").Value = 35617.52180241271
$ws.Range("AI13").Value = -227081.284486152
$ws.Range("AJ13").Value = -12.52138795235413
$ws.Range("BA13").Value = 0.9331624845914369
$ws.Range("BB13").Value = -14546.80934947369
$ws.Range("BG13").Value = 20.82149048541144
$ws.Range("BH13").Value = 20.34884085465865

$ws.Range("AH14").Value = 36981.55298073561
$ws.Range("AI14").Value = -212534.4751366783
$ws.Range("AJ14").Value = -12.52138795235413
$ws.Range("BA14").Value = 0.9477745788647102
$ws.Range("BB14").Value = -22079.09239733252
$ws.Range("BG14").Value = 20.34884085465865
$ws.Range("BH14").Value = 20.37712354942133

$ws.Range("AH15").Value = 30703.13424109436
$ws.Range("AI15").Value = -190455.3827393458
$ws.Range("AJ15").Value = -12.52138795235413
$ws.Range("BA15").Value = 0.9736183493953846
$ws.Range("BB15").Value = -25199.09552290288
$ws.Range("BG15").Value = 20.37712354942133
$ws.Range("BH15").Value = 20.96485846113388

$ws.Range("AH16").Value = 26032.49947366162
$ws.Range("AI16").Value = -165256.2872164429
$ws.Range("AJ16").Value = -12.52138795235413
$ws.Range("BA16").Value = 1.011777104389915
$ws.Range("BB16").Value = -26778.41970269778
$ws.Range("BG16").Value = 20.96485846113388
$ws.Range("BH16").Value = 22.03300056412582

$ws.Range("AH17").Value = 19744.03263057223
$ws.Range("AI17").Value = -138477.8675137452
$ws.Range("AJ17").Value = -12.52138795235413
$ws.Range("BA17").Value = 1.061317485883823
$ws.Range("BB17").Value = -23900.02137026326
$ws.Range("BG17").Value = 22.03300056412582
$ws.Range("BH17").Value = 23.34954780239644

$ws.Range("AH18").Value = 18726.97334774573
$ws.Range("AI18").Value = -114577.8461434819
$ws.Range("AJ18").Value = -12.52138795235413
$ws.Range("BA18").Value = 1.12336742583117
$ws.Range("BB18").Value = -25610.0576362883
$ws.Range("BG18").Value = 23.34954780239644
$ws.Range("BH18").Value = 25.26571607735314

$ws.Range("AH19").Value = 17305.62432390504
$ws.Range("AI19").Value = -88967.78850719359
$ws.Range("AJ19").Value = -12.52138795235413
$ws.Range("BA19").Value = 1.207839240417992
$ws.Range("BB19").Value = -26266.89058348348
$ws.Range("BG19").Value = 25.26571607735314
$ws.Range("BH19").Value = 27.86493428664652

$ws.Range("AH20").Value = 15339.6100186595
$ws.Range("AI20").Value = -62700.8979237101
$ws.Range("AJ20").Value = -12.52138795235413
$ws.Range("BA20").Value = 1.326632289745931
$ws.Range("BB20").Value = -25398.75199663377
$ws.Range("BG20").Value = 27.86493428664652
$ws.Range("BH20").Value = 31.2717180758355

$ws.Range("AH21").Value = 12587.38274117362
$ws.Range("AI21").Value = -37302.14592707633
$ws.Range("AJ21").Value = -12.52138795235413
$ws.Range("BA21").Value = 1.50350259473588
$ws.Range("BB21").Value = -22198.0169999653
$ws.Range("BG21").Value = 31.2717180758355
$ws.Range("BH21").Value = 35.68076560686083

$ws.Range("AH22").Value = 8614.414358475264
$ws.Range("AI22").Value = -15104.12892711103
$ws.Range("AJ22").Value = -12.52138795235413
$ws.Range("BA22").Value = 1.794330196609224
$ws.Range("BB22").Value = -15104.12892711105
$ws.Range("BG22").Value = 35.68076560686083
